{"js": "// Lattice-multiplication worksheet refresh: swap each cell's five lines\n// (top \"A x B\" expression, the split second-factor digits, the \"----\"\n// rule, and the two split first-factor digit rows) for a new fact, while\n// leaving the table's shape (5 rows x 3 columns) and per-run formatting\n// (sz=32) untouched. \\u000b below is the Word \"line break\" character that\n// backs each <w:br/> inside the single run/paragraph per cell.\nconst BR = \"\\u000b\";\n\n// New content for every cell, in row-major order (row0: col0,col1,col2;\n// row1: col0,col1,col2; ...). Each entry is the 5 lines for that cell.\nconst newCellLines = [\n  [\"18 x 53\", \"  5    3\", \"  ----\", \"1|    |\", \"8|    |\"],\n  [\"26 x 27\", \"  2    7\", \"  ----\", \"2|    |\", \"6|    |\"],\n  [\"92 x 62\", \"  6    2\", \"  ----\", \"9|    |\", \"2|    |\"],\n\n  [\"63 x 15\", \"  1    5\", \"  ----\", \"6|    |\", \"3|    |\"],\n  [\"48 x 43\", \"  4    3\", \"  ----\", \"4|    |\", \"8|    |\"],\n  [\"64 x 11\", \"  1    1\", \"  ----\", \"6|    |\", \"4|    |\"],\n\n  [\"68 x 93\", \"  9    3\", \"  ----\", \"6|    |\", \"8|    |\"],\n  [\"93 x 42\", \"  4    2\", \"  ----\", \"9|    |\", \"3|    |\"],\n  [\"57 x 45\", \"  4    5\", \"  ----\", \"5|    |\", \"7|    |\"],\n\n  [\"97 x 39\", \"  3    9\", \"  ----\", \"9|    |\", \"7|    |\"],\n  [\"57 x 53\", \"  5    3\", \"  ----\", \"5|    |\", \"7|    |\"],\n  [\"65 x 85\", \"  8    5\", \"  ----\", \"6|    |\", \"5|    |\"],\n\n  [\"66 x 71\", \"  7    1\", \"  ----\", \"6|    |\", \"6|    |\"],\n  [\"52 x 43\", \"  4    3\", \"  ----\", \"5|    |\", \"2|    |\"],\n  [\"90 x 82\", \"  8    2\", \"  ----\", \"9|    |\", \"0|    |\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n// `columnCount` isn't populated by this host; derive the grid shape from\n// the current values matrix instead (still 5 rows x 3 columns here).\ntable.load(\"values\");\nawait context.sync();\n\nconst rowCount = table.values.length;\nconst columnCount = rowCount > 0 ? table.values[0].length : 0;\n\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    if (idx >= newCellLines.length) break;\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    const text = newCellLines[idx].join(BR);\n    para.insertText(text, Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Lattice-multiplication worksheet refresh: swap each cell's five lines\n# (top \"A x B\" expression, the split second-factor digits, the \"----\"\n# rule, and the two split first-factor digit rows) for a new fact, while\n# leaving the table's shape (5 rows x 3 columns) and per-run formatting\n# (sz=32) untouched. [char]11 is the Word \"line break\" character that\n# backs each <w:br/> inside the single run/paragraph per cell.\n$vt = [char]11\n\n# New content for every cell, in row-major order (row0: col0,col1,col2;\n# row1: col0,col1,col2; ...). Each entry is the 5 lines for that cell.\n$newCellLines = @(\n    @(\"18 x 53\", \"  5    3\", \"  ----\", \"1|    |\", \"8|    |\"),\n    @(\"26 x 27\", \"  2    7\", \"  ----\", \"2|    |\", \"6|    |\"),\n    @(\"92 x 62\", \"  6    2\", \"  ----\", \"9|    |\", \"2|    |\"),\n\n    @(\"63 x 15\", \"  1    5\", \"  ----\", \"6|    |\", \"3|    |\"),\n    @(\"48 x 43\", \"  4    3\", \"  ----\", \"4|    |\", \"8|    |\"),\n    @(\"64 x 11\", \"  1    1\", \"  ----\", \"6|    |\", \"4|    |\"),\n\n    @(\"68 x 93\", \"  9    3\", \"  ----\", \"6|    |\", \"8|    |\"),\n    @(\"93 x 42\", \"  4    2\", \"  ----\", \"9|    |\", \"3|    |\"),\n    @(\"57 x 45\", \"  4    5\", \"  ----\", \"5|    |\", \"7|    |\"),\n\n    @(\"97 x 39\", \"  3    9\", \"  ----\", \"9|    |\", \"7|    |\"),\n    @(\"57 x 53\", \"  5    3\", \"  ----\", \"5|    |\", \"7|    |\"),\n    @(\"65 x 85\", \"  8    5\", \"  ----\", \"6|    |\", \"5|    |\"),\n\n    @(\"66 x 71\", \"  7    1\", \"  ----\", \"6|    |\", \"6|    |\"),\n    @(\"52 x 43\", \"  4    3\", \"  ----\", \"5|    |\", \"2|    |\"),\n    @(\"90 x 82\", \"  8    2\", \"  ----\", \"9|    |\", \"0|    |\")\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$idx = 0\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    for ($c = 1; $c -le $table.Columns.Count; $c++) {\n        if ($idx -ge $newCellLines.Count) { break }\n        $lines = $newCellLines[$idx]\n        $newText = [string]::Join($vt, $lines)\n        $table.Cell($r, $c).Range.Text = $newText\n        $idx++\n    }\n}\n\nWrite-Output \"Updated $idx cells\"\n"}
